$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("templates")

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 'leaves viewers out in the cold and undermines some phenomenal performances .'
$ws.Range("D2").Value = 'leaves viewers out in the cold and {mask} some {mask} performances .'
$ws.Range("E2").Value = 'leaves viewers out in the cold and {neg_verb} some {pos_adj} performances .'

$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 'one gets the impression the creators of don''t ask don''t tell laughed a hell of a lot at their own jokes . too bad none of it is funny .'
$ws.Range("D3").Value = 'one {mask} the impression the creators of do n''t ask do n''t tell laughed a hell of a lot at their own jokes . too {mask} none of it is funny .'
$ws.Range("E3").Value = 'one {pos_verb} the impression the creators of do n''t ask do n''t tell laughed a hell of a lot at their own jokes . too {neg_adj} none of it is funny .'

$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 'connoisseurs of chinese film will be pleased to discover that tian''s meticulous talent has not withered during his enforced hiatus .'
$ws.Range("D4").Value = 'connoisseurs of chinese film will be pleased to {mask} that tian ''s meticulous talent has not {mask} during his enforced hiatus .'
$ws.Range("E4").Value = 'connoisseurs of chinese film will be pleased to {pos_verb} that tian ''s meticulous talent has not {neg_verb} during his enforced hiatus .'

$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 'a moving and not infrequently breathtaking film .'
$ws.Range("D5").Value = 'a {mask} and not infrequently {mask} film .'
$ws.Range("E5").Value = 'a {pos_adj} and not infrequently {pos_adj} film .'

$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 'caine makes us watch as his character awakens to the notion that to be human is eventually to have to choose .'
$ws.Range("D6").Value = 'caine makes us watch as his character {mask} to the notion that to be human is eventually to {mask} to choose .'
$ws.Range("E6").Value = 'caine makes us watch as his character {pos_verb} to the notion that to be human is eventually to {neg_verb} to choose .'

$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 'nothing but an episode of smackdown ! in period costume and with a bigger budget .'
$ws.Range("D7").Value = 'nothing but an episode of smackdown ! in period costume and with a {mask} budget .'
$ws.Range("E7").Value = 'nothing but an episode of smackdown ! in period costume and with a {pos_adj} budget .'

$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 'yes they can swim , the title is merely anne-sophie birot''s off-handed way of saying girls find adolescence difficult to wade through .'
$ws.Range("D8").Value = 'yes they can swim , the title is merely anne - sophie birot ''s {mask} - {mask} way of saying girls find adolescence difficult to wade through .'
$ws.Range("E8").Value = 'yes they can swim , the title is merely anne - sophie birot ''s {neg_adj} - {neg_adj} way of saying girls find adolescence difficult to wade through .'

$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 'a rip-off twice removed , modeled after [seagal''s] earlier copycat under siege , sometimes referred to as die hard on a boat .'
$ws.Range("D9").Value = 'a rip - off twice {mask} , {mask} after [ seagal ''s ] earlier copycat under siege , sometimes referred to as die hard on a boat .'
$ws.Range("E9").Value = 'a rip - off twice {neg_verb} , {neg_verb} after [ seagal ''s ] earlier copycat under siege , sometimes referred to as die hard on a boat .'

$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 'jolie gives it that extra little something that makes it worth checking out at theaters , especially if you''re in the mood for something more comfortable than challenging .'
$ws.Range("D10").Value = 'jolie {mask} it that extra little something that makes it worth checking out at theaters , especially if you ''re in the mood for something more comfortable than {mask} .'
$ws.Range("E10").Value = 'jolie {pos_verb} it that extra little something that makes it worth checking out at theaters , especially if you ''re in the mood for something more comfortable than {pos_adj} .'

$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 'a dreadful day in irish history is given passionate , if somewhat flawed , treatment .'
$ws.Range("D11").Value = 'a dreadful day in irish history is {mask} {mask} , if somewhat flawed , treatment .'
$ws.Range("E11").Value = 'a dreadful day in irish history is {neg_verb} {pos_adj} , if somewhat flawed , treatment .'

$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 'what might have been readily dismissed as the tiresome rant of an aging filmmaker still thumbing his nose at convention takes a surprising , subtle turn at the midway point .'
$ws.Range("D12").Value = 'what might have been readily {mask} as the {mask} rant of an aging filmmaker still thumbing his nose at convention takes a surprising , subtle turn at the midway point .'
$ws.Range("E12").Value = 'what might have been readily {neg_verb} as the {neg_adj} rant of an aging filmmaker still thumbing his nose at convention takes a surprising , subtle turn at the midway point .'

$ws.Range("B13").Value = 0
$ws.Range("C13").Value = ' one look at a girl in tight pants and big tits and you turn stupid ?  um? . . isn''t that the basis for the entire plot ?'
$ws.Range("D13").Value = ' one look at a girl in tight pants and big tits and you turn {mask} ?  um ? . . is n''t that the basis for the {mask} plot ?'
$ws.Range("E13").Value = ' one look at a girl in tight pants and big tits and you turn {neg_adj} ?  um ? . . is n''t that the basis for the {neg_adj} plot ?'

$ws.Range("B14").Value = 1
$ws.Range("C14").Value = 'the stunt work is top-notch ; the dialogue and drama often food-spittingly funny .'
$ws.Range("D14").Value = 'the stunt work is top - {mask} ; the dialogue and drama often food - spittingly {mask} .'
$ws.Range("E14").Value = 'the stunt work is top - {pos_adj} ; the dialogue and drama often food - spittingly {pos_adj} .'

$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 'a film of delicate interpersonal dances .'
$ws.Range("D15").Value = 'a film of {mask} {mask} dances .'
$ws.Range("E15").Value = 'a film of {pos_adj} {pos_adj} dances .'

$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 'the rules of attraction gets us too drunk on the party favors to sober us up with the transparent attempts at moralizing .'
$ws.Range("D16").Value = 'the rules of attraction {mask} us too {mask} on the party favors to sober us up with the transparent attempts at moralizing .'
$ws.Range("E16").Value = 'the rules of attraction {pos_verb} us too {neg_adj} on the party favors to sober us up with the transparent attempts at moralizing .'

$ws.Range("B17").Value = 1
$ws.Range("C17").Value = 'twohy knows how to inflate the mundane into the scarifying , and gets full mileage out of the rolling of a stray barrel or the unexpected blast of a phonograph record .'
$ws.Range("D17").Value = 'twohy {mask} how to inflate the mundane into the scarifying , and {mask} full mileage out of the rolling of a stray barrel or the unexpected blast of a phonograph record .'
$ws.Range("E17").Value = 'twohy {pos_verb} how to inflate the mundane into the scarifying , and {pos_verb} full mileage out of the rolling of a stray barrel or the unexpected blast of a phonograph record .'

$ws.Range("B18").Value = 0
$ws.Range("C18").Value = 'this new zealand coming-of-age movie isn''t really about anything .'
$ws.Range("D18").Value = 'this new zealand {mask} - of - age movie is n''t really about anything .'
$ws.Range("E18").Value = 'this new zealand {pos_verb} - of - age movie is n''t really about anything .'

$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 'death to smoochy tells a moldy-oldie , not-nearly -as-nasty -as-it- thinks-it-is joke . over and over again .'
$ws.Range("D19").Value = 'death to smoochy {mask} a moldy - oldie , not - nearly -as - {mask} -as - it- thinks - it - is joke . over and over again .'
$ws.Range("E19").Value = 'death to smoochy {neg_verb} a moldy - oldie , not - nearly -as - {neg_adj} -as - it- thinks - it - is joke . over and over again .'
